# Reparando boton cancelar de confpos.html
# Populates the "ruta_actual" sheet with the first stop of the newly
# started route and records that route in "rutas_registros".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ruta_actual": header info (B1/C1) + first client row (row 3)
# ---------------------------------------------------------------------
$wsRuta = $wb.Worksheets.Item("ruta_actual")

# Cells that hold digit-only strings need to stay text, not be coerced
# into numbers, so force a text number format before assigning them.
$wsRuta.Range("B1").NumberFormat = "@"
$wsRuta.Range("B1").Value = "20240816"
$wsRuta.Range("C1").NumberFormat = "@"
$wsRuta.Range("C1").Value = "ruta TEST"

$wsRuta.Range("A3").NumberFormat = "@"
$wsRuta.Range("A3").Value = "20240816"
$wsRuta.Range("B3").Value = 1
$wsRuta.Range("C3").NumberFormat = "@"
$wsRuta.Range("C3").Value = "16.742.249-7"
$wsRuta.Range("D3").Value = "Isaias Beroiza Mora"
$wsRuta.Range("E3").Value = "colaco sn km3 parcela 9"
$wsRuta.Range("F3").Value = "Calbuco"
$wsRuta.Range("G3").NumberFormat = "@"
$wsRuta.Range("G3").Value = "88809703"
$wsRuta.Range("H3").Value = "por buscar"
$wsRuta.Range("I3").Value = "ok"
$wsRuta.Range("J3").NumberFormat = "@"
$wsRuta.Range("J3").Value = "20240813"
$wsRuta.Range("K3").Value = 120
$wsRuta.Range("L3").NumberFormat = "@"
$wsRuta.Range("L3").Value = "2024-12-01"

# ---------------------------------------------------------------------
# Sheet "rutas_registros": new route entry for 20240816 / ruta TEST
# ---------------------------------------------------------------------
$wsRegistros = $wb.Worksheets.Item("rutas_registros")

$wsRegistros.Range("A5").NumberFormat = "@"
$wsRegistros.Range("A5").Value = "20240816"
$wsRegistros.Range("B5").Value = "ruta TEST"

Write-Output "edit applied"
